$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("D2").Value = "297.39"
$ws.Range("E2").Value = "-1.11%"
$ws.Range("G2").Value = "17"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("D3").Value = "31.68"
$ws.Range("E3").Value = "0.75%"
$ws.Range("G3").Value = "17"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("D4").Value = "5.075"
$ws.Range("E4").Value = "-1.64%"
$ws.Range("G4").Value = "17"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08030"
$ws.Range("E5").Value = "8.76%"
$ws.Range("G5").Value = "17"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("D6").Value = "2.620"
$ws.Range("E6").Value = "23.26%"
$ws.Range("G6").Value = "17"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("D7").Value = "7.793"
$ws.Range("E7").Value = "-1.34%"
$ws.Range("G7").Value = "17"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("D8").Value = "3.817"
$ws.Range("E8").Value = "1.28%"
$ws.Range("G8").Value = "17"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9243"
$ws.Range("E9").Value = "-0.25%"
$ws.Range("G9").Value = "17"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1755"
$ws.Range("E10").Value = "2.63%"
$ws.Range("G10").Value = "17"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07370"
$ws.Range("E11").Value = "-3.51%"
$ws.Range("G11").Value = "17"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08910"
$ws.Range("E12").Value = "9.45%"
$ws.Range("G12").Value = "17"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03062"
$ws.Range("E13").Value = "1.47%"
$ws.Range("G13").Value = "17"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("E14").Value = "0.85%"
$ws.Range("G14").Value = "17"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001493"
$ws.Range("E15").Value = "-1.35%"
$ws.Range("G15").Value = "17"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006054"
$ws.Range("E16").Value = "-0.53%"
$ws.Range("G16").Value = "17"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("D17").Value = "3.551"
$ws.Range("E17").Value = "2.67%"
$ws.Range("G17").Value = "17"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("E18").Value = "1.15%"
$ws.Range("G18").Value = "17"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("E19").Value = "0.79%"
$ws.Range("G19").Value = "17"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1346"
$ws.Range("E20").Value = "-0.10%"
$ws.Range("G20").Value = "17"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("D21").Value = "4.000"
$ws.Range("E21").Value = "-14.04%"
$ws.Range("G21").Value = "17"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1646"
$ws.Range("E22").Value = "3.93%"
$ws.Range("G22").Value = "17"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04598"
$ws.Range("E23").Value = "-1.19%"
$ws.Range("G23").Value = "17"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001243"
$ws.Range("E24").Value = "2.00%"
$ws.Range("G24").Value = "17"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004436"
$ws.Range("G25").Value = "17"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.80%"
$ws.Range("G26").Value = "17"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003408"
$ws.Range("E27").Value = "82.07%"
$ws.Range("G27").Value = "17"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "17"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "17"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "17"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "17"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "17"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "17"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "17"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "17"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "17"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "17"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "17"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01767"
$ws.Range("E39").Value = "1.92%"
$ws.Range("G39").Value = "17"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04504"
$ws.Range("E40").Value = "-0.32%"
$ws.Range("G40").Value = "17"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006757"
$ws.Range("E41").Value = "-4.76%"
$ws.Range("G41").Value = "17"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1348"
$ws.Range("E42").Value = "0.47%"
$ws.Range("G42").Value = "17"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.33%"
$ws.Range("G43").Value = "17"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009845"
$ws.Range("E44").Value = "-6.48%"
$ws.Range("G44").Value = "17"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006458"
$ws.Range("E45").Value = "3.00%"
$ws.Range("G45").Value = "17"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").Value = "-0.15%"
$ws.Range("G46").Value = "17"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "0.008737"
$ws.Range("E47").Value = "24.98%"
$ws.Range("G47").Value = "17"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "0.8206"
$ws.Range("E48").Value = "-57.45%"
$ws.Range("G48").Value = "17"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002097"
$ws.Range("E49").Value = "-0.15%"
$ws.Range("G49").Value = "17"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001997"
$ws.Range("E50").Value = "-0.08%"
$ws.Range("G50").Value = "17"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "17"
